$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows whose "Sending cluster" is MuSCs (original rows 6-9).
# Remaining rows (originally 2-5, FAPs sending cluster) keep their row numbers.
$ws.Rows("6:9").Delete()

# Updated TPM-derived numeric values for the four remaining data rows.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1.594873333333333
$ws.Range("N2").Value = 4.78462
$ws.Range("O2").Value = 0.09372679355272211
$ws.Range("P2").Value = 0.09372679355272213
$ws.Range("Q2").Value = 2.449360214006667
$ws.Range("R2").Value = 22.04424192606
$ws.Range("S2").Value = 0.09372679355272211
$ws.Range("T2").Value = 0.09372679355272213

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 4.578777
$ws.Range("N3").Value = 13.736331
$ws.Range("O3").Value = 0.2690834924840127
$ws.Range("P3").Value = 0.2690834924840128
$ws.Range("Q3").Value = 7.031952932067001
$ws.Range("R3").Value = 63.28757638860301
$ws.Range("S3").Value = 0.2690834924840127
$ws.Range("T3").Value = 0.2690834924840128

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 4.495828
$ws.Range("N4").Value = 13.487484
$ws.Range("O4").Value = 0.2642087832291055
$ws.Range("P4").Value = 0.2642087832291055
$ws.Range("Q4").Value = 6.904562263388001
$ws.Range("R4").Value = 62.14106037049201
$ws.Range("S4").Value = 0.2642087832291055
$ws.Range("T4").Value = 0.2642087832291055

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 6.346716
$ws.Range("N5").Value = 19.040148
$ws.Range("O5").Value = 0.3729809307341596
$ws.Range("P5").Value = 0.3729809307341597
$ws.Range("Q5").Value = 9.747102378036001
$ws.Range("R5").Value = 87.723921402324
$ws.Range("S5").Value = 0.3729809307341596
$ws.Range("T5").Value = 0.3729809307341597
